# Applies the "Optuna Attempt (go back with original)" edits:
#  - Updates Inventory Coverage (H) and Seasonality Index (L) values on the
#    "Forecast Comparison" sheet for weeks W8..W23 (rows 2-17).
#  - Updates the "Total Forecast (16 Weeks)" / "Total Forecast (8 Weeks)"
#    figures on the "Summary" sheet (B9, B10), keeping them as text values.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: column H (Inventory Coverage) ---
$wsForecast.Range("H2").Value  = 33.33
$wsForecast.Range("H3").Value  = 20.42
$wsForecast.Range("H4").Value  = 14.76
$wsForecast.Range("H5").Value  = 11.1
$wsForecast.Range("H6").Value  = 10.1
$wsForecast.Range("H7").Value  = 11.28
$wsForecast.Range("H8").Value  = 13.53
$wsForecast.Range("H9").Value  = 9.52
$wsForecast.Range("H10").Value = 6.87
$wsForecast.Range("H11").Value = 5.87
$wsForecast.Range("H12").Value = 6.04
$wsForecast.Range("H13").Value = 6.63
$wsForecast.Range("H14").Value = 4.28
$wsForecast.Range("H15").Value = 4.56
$wsForecast.Range("H16").Value = 4.27
$wsForecast.Range("H17").Value = 2.58

# --- Forecast Comparison sheet: column L (Seasonality Index) ---
$wsForecast.Range("L2").Value  = 0.84
$wsForecast.Range("L3").Value  = 1.12
$wsForecast.Range("L4").Value  = 1.03
$wsForecast.Range("L5").Value  = 0.98
$wsForecast.Range("L6").Value  = 1.14
$wsForecast.Range("L7").Value  = 1.18
$wsForecast.Range("L9").Value  = 1.16
$wsForecast.Range("L10").Value = 1.02
$wsForecast.Range("L11").Value = 0.88
$wsForecast.Range("L12").Value = 1.13
$wsForecast.Range("L13").Value = 1.15
$wsForecast.Range("L14").Value = 1.08
$wsForecast.Range("L15").Value = 0.96
$wsForecast.Range("L16").Value = 0.91
$wsForecast.Range("L17").Value = 1.17

# --- Summary sheet: keep these as text values (not numbers) ---
$wsSummary.Range("B9").NumberFormat  = "@"
$wsSummary.Range("B9").Value  = "7"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "4"

Write-Host "Applied Optuna attempt edits to Forecast Comparison and Summary sheets."
